$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 1-5 (L_D, L, D, CL, CD) with refined values
$ws.Range("B1").Value = 4.004048952012076
$ws.Range("C1").Value = 3.383645344629007
$ws.Range("D1").Value = 3.624549554035144

$ws.Range("B2").Value = 3223050.483304044
$ws.Range("C2").Value = 3208196.205433025
$ws.Range("D2").Value = 3214446.834770465

$ws.Range("B3").Value = 804947.8220500844
$ws.Range("C3").Value = 948147.8933734945
$ws.Range("D3").Value = 886854.155764509

$ws.Range("B4").Value = 0.09901473025964723
$ws.Range("C4").Value = 0.09855839477119624
$ws.Range("D4").Value = 0.0987504191844053

$ws.Range("B5").Value = 0.02472865128431841
$ws.Range("C5").Value = 0.02912787385582293
$ws.Range("D5").Value = 0.0272448804223047

# Row 6 (S) unchanged

# Row 7 (V) updated
$ws.Range("B7").Value = 10422.67289148332
$ws.Range("C7").Value = 10422.67289148332
$ws.Range("D7").Value = 10422.67289148332

# Row 8 (volEff) updated
$ws.Range("B8").Value = 0.07999997241986209
$ws.Range("C8").Value = 0.07999997241986209
$ws.Range("D8").Value = 0.07999997241986209

# Rows 9-11 (s, l, s_l) unchanged

# New rows 12-15: Cp_base, P_base, D_base, S_base
$ws.Range("A12").Value = "Cp_base"
$ws.Range("B12").Value = -0.04297219340437737
$ws.Range("C12").Value = -0.04297219340437737
$ws.Range("D12").Value = -0.04297219340437737

$ws.Range("A13").Value = "P_base"
$ws.Range("B13").Value = 179.4863027032873
$ws.Range("C13").Value = 179.4863027032873
$ws.Range("D13").Value = 179.4863027032873

$ws.Range("A14").Value = "D_base"
$ws.Range("B14").Value = 36824.8054614949
$ws.Range("C14").Value = 36824.8054614949
$ws.Range("D14").Value = 36824.8054614949

$ws.Range("A15").Value = "S_base"
$ws.Range("B15").Value = 205.1677755175045
$ws.Range("C15").Value = 205.1677755175045
$ws.Range("D15").Value = 205.1677755175045
